$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Total Sales" text for rows 3 and 4 to reflect corrected values
$ws.Range("F3").Value = "2*1000=2000"
$ws.Range("F4").Value = "2*2000=4000"

# Update the selected cell to F4 (last active cell after edit)
$ws.Range("F4").Select()
